$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 367. This shifts the existing row 367
# (and everything below it) down by one, turning old row 367 into row 368,
# old row 431 into row 432, etc., and growing the used range to R432.
$ws.Rows.Item(367).Insert()

# The new row 367 is a new weekly observation for the same market/product
# series as the (now shifted) row 368, carrying the same attributes except
# for the date (column D) and volume (column J).
$ws.Range("A367").Value2 = $ws.Range("A368").Value2
$ws.Range("B367").Value2 = $ws.Range("B368").Value2
$ws.Range("C367").Value2 = $ws.Range("C368").Value2
$ws.Range("D367").Value2 = 45258
$ws.Range("D367").NumberFormat = $ws.Range("D368").NumberFormat
$ws.Range("E367").Value2 = $ws.Range("E368").Value2
$ws.Range("F367").Value2 = $ws.Range("F368").Value2
$ws.Range("G367").Value2 = $ws.Range("G368").Value2
$ws.Range("H367").Value2 = $ws.Range("H368").Value2
$ws.Range("I367").Value2 = $ws.Range("I368").Value2
$ws.Range("J367").Value2 = 250
$ws.Range("K367").Value2 = $ws.Range("K368").Value2
$ws.Range("L367").Value2 = $ws.Range("L368").Value2
$ws.Range("M367").Value2 = $ws.Range("M368").Value2
$ws.Range("N367").Value2 = $ws.Range("N368").Value2
$ws.Range("O367").Value2 = $ws.Range("O368").Value2
$ws.Range("P367").Value2 = $ws.Range("P368").Value2
$ws.Range("Q367").Value2 = $ws.Range("Q368").Value2
$ws.Range("R367").Value2 = $ws.Range("R368").Value2
